$wb = $excel.ActiveWorkbook

# ALC row 28
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 4964.577
$ws.Range("I28").Value = 762.2308
$ws.Range("J28").Value = 9166.923000000001
$ws.Range("K28").Value = 762.2308
$ws.Range("L28").Value = 9166.923000000001
$ws.Range("M28").Value = -277.2308
$ws.Range("N28").Value = -10136.923

# ALC row 31
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 2021.1428
$ws.Range("I31").Value = 2021.1428
$ws.Range("K31").Value = 6063.428400000001
$ws.Range("M31").Value = -5833.428400000001

# ALC row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 7250460
$ws.Range("I62").Value = 7940728
$ws.Range("J62").Value = 2645.5
$ws.Range("K62").Value = 7940728
$ws.Range("L62").Value = 2645.5
$ws.Range("M62").Value = -7940104
$ws.Range("N62").Value = -3893.5

# ALC row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 7250460
$ws.Range("I65").Value = 7940728
$ws.Range("J65").Value = 2645.5
$ws.Range("K65").Value = 39703640
$ws.Range("L65").Value = 13227.5
$ws.Range("M65").Value = -39700520
$ws.Range("N65").Value = -19467.5

# ALC row 107
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 888.9211
$ws.Range("I107").Value = 1004.2143
$ws.Range("K107").Value = 1004.2143
$ws.Range("M107").Value = 915.7857

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 5558719
$ws.Range("I116").Value = 10103636
$ws.Range("J116").Value = 3820.2222
$ws.Range("K116").Value = 10103636
$ws.Range("L116").Value = 3820.2222
$ws.Range("M116").Value = -10100194
$ws.Range("N116").Value = -10704.2222

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 7048.457
$ws.Range("I137").Value = 4573.2856
$ws.Range("J137").Value = 16949.143
$ws.Range("K137").Value = 13719.8568
$ws.Range("L137").Value = 50847.429
$ws.Range("M137").Value = -11169.8568
$ws.Range("N137").Value = -55947.429

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3029.726
$ws.Range("I138").Value = 979.12823
$ws.Range("J138").Value = 5381.8823
$ws.Range("K138").Value = 2937.38469
$ws.Range("L138").Value = 16145.6469
$ws.Range("M138").Value = 2202.61531
$ws.Range("N138").Value = -26425.6469

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 1053.6451
$ws.Range("I141").Value = 1089.0714
$ws.Range("K141").Value = 3267.2142
$ws.Range("M141").Value = 1912.7858

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7548392
$ws.Range("I32").Value = 7044292
$ws.Range("J32").Value = 11127500
$ws.Range("K32").Value = 7044292
$ws.Range("L32").Value = 11127500
$ws.Range("M32").Value = -7044005
$ws.Range("N32").Value = -11128074

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2933.3333
$ws.Range("I45").Value = 2925
$ws.Range("J45").Value = 3000
$ws.Range("K45").Value = 2925
$ws.Range("L45").Value = 3000
$ws.Range("M45").Value = -2548
$ws.Range("N45").Value = -3754

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3849.4
$ws.Range("I74").Value = 3873.182
$ws.Range("J74").Value = 3784
$ws.Range("K74").Value = 3873.182
$ws.Range("L74").Value = 3784
$ws.Range("M74").Value = -2999.182
$ws.Range("N74").Value = -5532

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 3849.4
$ws.Range("I77").Value = 3873.182
$ws.Range("J77").Value = 3784
$ws.Range("K77").Value = 19365.91
$ws.Range("L77").Value = 18920
$ws.Range("M77").Value = -14997.91
$ws.Range("N77").Value = -27656

# ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 4908.8
$ws.Range("I110").Value = 2492.125
$ws.Range("J110").Value = 5787.591
$ws.Range("K110").Value = 2492.125
$ws.Range("L110").Value = 5787.591
$ws.Range("M110").Value = -447.125
$ws.Range("N110").Value = -9877.591

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 540569
$ws.Range("I132").Value = 604310.7
$ws.Range("K132").Value = 1812932.1
$ws.Range("M132").Value = -1810402.1

# BSM row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 9525248
$ws.Range("I20").Value = 15874247
$ws.Range("J20").Value = 1748.6666
$ws.Range("K20").Value = 15874247
$ws.Range("L20").Value = 1748.6666
$ws.Range("M20").Value = -15874000
$ws.Range("N20").Value = -2242.6666

# BSM row 80
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 16680366
$ws.Range("J80").Value = 25660918
$ws.Range("L80").Value = 25660918
$ws.Range("N80").Value = -25662914

# BSM row 83
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 16680366
$ws.Range("J83").Value = 25660918
$ws.Range("L83").Value = 128304590
$ws.Range("N83").Value = -128314574

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1168804.6
$ws.Range("I134").Value = 1224917.1
$ws.Range("J134").Value = 18500
$ws.Range("K134").Value = 3674751.3
$ws.Range("L134").Value = 55500
$ws.Range("M134").Value = -3672216.3
$ws.Range("N134").Value = -60570

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 21282244
$ws.Range("I58").Value = 27781638
$ws.Range("J58").Value = 11500.546
$ws.Range("K58").Value = 27781638
$ws.Range("L58").Value = 11500.546
$ws.Range("M58").Value = -27781435
$ws.Range("N58").Value = -11906.546

# CRP row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1470.6875
$ws.Range("I107").Value = 725.25
$ws.Range("J107").Value = 2216.125
$ws.Range("K107").Value = 725.25
$ws.Range("L107").Value = 2216.125
$ws.Range("M107").Value = 1194.75
$ws.Range("N107").Value = -6056.125

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 5488.522
$ws.Range("I132").Value = 4546.737
$ws.Range("J132").Value = 9962
$ws.Range("K132").Value = 13640.211
$ws.Range("L132").Value = 29886
$ws.Range("M132").Value = -11110.211
$ws.Range("N132").Value = -34946

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 21282244
$ws.Range("I136").Value = 27781638
$ws.Range("J136").Value = 11500.546
$ws.Range("K136").Value = 83344914
$ws.Range("L136").Value = 34501.638
$ws.Range("M136").Value = -83342364
$ws.Range("N136").Value = -39601.638

# CUL row 6
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 804.44446
$ws.Range("I6").Value = 804.44446
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 2413.33338
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -2300.33338
$ws.Range("N6").ClearContents()

# CUL row 38
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 64.5
$ws.Range("J38").Value = 64.5
$ws.Range("L38").Value = 193.5
$ws.Range("N38").Value = -887.5

# CUL row 39
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 6850
$ws.Range("J39").Value = 6850
$ws.Range("L39").Value = 20550
$ws.Range("N39").Value = -21138

# CUL row 50
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 1056.7142
$ws.Range("I50").Value = 1056.7142
$ws.Range("K50").Value = 3170.1426
$ws.Range("M50").Value = -2689.1426

# CUL row 53
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H53").Value = 1056.7142
$ws.Range("I53").Value = 1056.7142
$ws.Range("K53").Value = 3170.1426
$ws.Range("M53").Value = -2689.1426

# CUL row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2502500
$ws.Range("J68").Value = 2502500
$ws.Range("L68").Value = 7507500
$ws.Range("N68").Value = -7509122

# CUL row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 2502500
$ws.Range("J71").Value = 2502500
$ws.Range("L71").Value = 22522500
$ws.Range("N71").Value = -22530612

# CUL row 74
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 27499.25
$ws.Range("J74").Value = 27499.25
$ws.Range("L74").Value = 82497.75
$ws.Range("N74").Value = -84619.75

# CUL row 77
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H77").Value = 27499.25
$ws.Range("J77").Value = 27499.25
$ws.Range("L77").Value = 247493.25
$ws.Range("N77").Value = -258101.25

# CUL row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2988.1875
$ws.Range("I132").Value = 1960.6666
$ws.Range("J132").Value = 3225.3076
$ws.Range("K132").Value = 17645.9994
$ws.Range("L132").Value = 29027.7684
$ws.Range("M132").Value = -15115.9994
$ws.Range("N132").Value = -34087.7684

# GSM row 2
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 78.78570999999999
$ws.Range("I2").Value = 77.3
$ws.Range("J2").Value = 82.5
$ws.Range("K2").Value = 77.3
$ws.Range("L2").Value = 82.5
$ws.Range("M2").Value = 35.7
$ws.Range("N2").Value = -308.5

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 14709015
$ws.Range("I132").Value = 20410988
$ws.Range("J132").Value = 3924.6316
$ws.Range("K132").Value = 61232964
$ws.Range("L132").Value = 11773.8948
$ws.Range("M132").Value = -61230434
$ws.Range("N132").Value = -16833.8948

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5914.0356
$ws.Range("I61").Value = 5129.273
$ws.Range("K61").Value = 5129.273
$ws.Range("M61").Value = -4927.273

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 5914.0356
$ws.Range("I113").Value = 5129.273
$ws.Range("K113").Value = 5129.273
$ws.Range("M113").Value = -2959.273

# LTW row 127
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2714.8254
$ws.Range("I132").Value = 2816.549
$ws.Range("J132").Value = 2282.5
$ws.Range("K132").Value = 8449.647000000001
$ws.Range("L132").Value = 6847.5
$ws.Range("M132").Value = -5919.647000000001
$ws.Range("N132").Value = -11907.5

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 28306800
$ws.Range("I136").Value = 12199543
$ws.Range("J136").Value = 83339930
$ws.Range("K136").Value = 36598629
$ws.Range("L136").Value = 250019790
$ws.Range("M136").Value = -36596079
$ws.Range("N136").Value = -250024890

# WVR row 95
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 30126.4
$ws.Range("J95").Value = 30126.4
$ws.Range("L95").Value = 30126.4
$ws.Range("N95").Value = -35618.4

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6342.7896
$ws.Range("I132").Value = 4519.6562
$ws.Range("K132").Value = 13558.9686
$ws.Range("M132").Value = -11028.9686

# WVR row 135
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 57600.5
$ws.Range("J135").Value = 57600.5
$ws.Range("L135").Value = 57600.5
$ws.Range("N135").Value = -67740.5

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 7046293
$ws.Range("I136").Value = 8475464
$ws.Range("J136").Value = 19534.5
$ws.Range("K136").Value = 25426392
$ws.Range("L136").Value = 58603.5
$ws.Range("M136").Value = -25423842
$ws.Range("N136").Value = -63703.5
